# The deck ships two DrawingML themes:
#   ppt/theme/theme1.xml -> "Office Theme" (clrScheme "Office", blue accents)
#   ppt/theme/theme2.xml -> "Integral"     (clrScheme "Red Violet", pink/purple accents)
# and theme2.xml is the one actually wired to the slide master / slides.
#
# The target edit swaps the two themes' content, i.e. the live theme (theme2.xml,
# reached here through ThemeColorScheme) ends up holding the "Office" blue color
# values instead of the "Red Violet" ones. fontScheme/fmtScheme are identical
# between the two themes already, so the only substantive change is the 12
# scheme colors (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink).

function Hex2RGB($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Target ("Office") theme colors, in ThemeColorScheme index order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink
$officeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$p = $ppt.ActivePresentation
$tcs = $p.Slides.Item(1).ThemeColorScheme

for ($i = 1; $i -le $officeColors.Count; $i++) {
    $tcs.Item($i).RGB = Hex2RGB $officeColors[$i - 1]
}
